# Corrected mistake in experiment.xlsx:
# The "cells" worksheet was missing a Gate Fraction value for the first
# data row (row 2, FCFiles/data_001.fcs). Fill in the missing value to
# match the other rows in the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cells")

$ws.Range("F2").Value = 0.3

# Reflect where the user's cursor ended up after making the correction.
$ws.Activate()
$ws.Range("E12").Select()
